# Updates cryptos list values (coin price/volume snapshot refresh)
# Applies the upstream diff cell-by-cell via the Excel COM object model.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.201.41"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'2.378.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.40%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'549.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.06%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'139.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.75%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.02%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.529"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.31%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.379.82"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.29%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +2.56%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +1.33%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'5.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.32%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.349"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.83%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'25.12"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.42%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +1.05%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'2.790.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.01%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'61.134.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.07%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.386.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.03%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'10.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.88%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").Value = "'4.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.73%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'321.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.88%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'6.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.54%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D24").Value = "'64.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.37%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -11.06%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'8.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +2.26%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -0.07%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.487.93"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.70%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'8.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +1.08%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'507.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -3.35%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.0₃0890"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -4.26%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +3.34%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -3.55%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -0.60%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -3.48%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +0.06%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'4.68"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.31%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").Value = "'5.41"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.94%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").Value = "'0.379"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.86%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'1.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.77%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'18.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +2.91%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'146.67"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +5.41%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D44").Value = "'41.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +2.90%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'148.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +5.67%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'3.61"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.35%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -2.98%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.0523"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.63%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'19.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -4.00%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.577"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.29%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +0.52%  "
$ws.Range("E51").Style = "Normal"

Write-Host "Applied cryptos update"
